$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '26.198.14'
$ws.Range('E2').Value = '  -4.02%  '
$ws.Range('D3').Value = '1.659.06'
$ws.Range('E3').Value = '  -2.72%  '
Set-TextValue $ws 'D4' '1.004'
$ws.Range('E4').Value = '  +0.09%  '
Set-TextValue $ws 'D5' '218.18'
$ws.Range('E5').Value = '  -2.56%  '
Set-TextValue $ws 'D6' '0.5158'
$ws.Range('E6').Value = '  -2.79%  '
$ws.Range('E7').Value = '  +0.11%  '
Set-TextValue $ws 'D8' '0.2588'
$ws.Range('E8').Value = '  -2.61%  '
Set-TextValue $ws 'D9' '0.06447'
$ws.Range('E9').Value = '  -1.95%  '
Set-TextValue $ws 'D10' '19.96'
$ws.Range('E10').Value = '  -3.76%  '
Set-TextValue $ws 'D11' '0.07800'
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('D12').Value = '1.658.65'
$ws.Range('E12').Value = '  -2.86%  '
Set-TextValue $ws 'D13' '4.296'
$ws.Range('E13').Value = '  -4.80%  '
$ws.Range('D14').Value = '1.886.80'
$ws.Range('E14').Value = '  -2.74%  '
Set-TextValue $ws 'D15' '0.5555'
$ws.Range('E15').Value = '  -3.68%  '
$ws.Range('D16').Value = '0.0₅8068'
$ws.Range('E16').Value = '  -0.86%  '
Set-TextValue $ws 'D17' '64.28'
$ws.Range('E17').Value = '  -4.86%  '
$ws.Range('D18').Value = '26.200.86'
$ws.Range('E18').Value = '  -4.04%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws 'D19' '211.42'
$ws.Range('E19').Value = '  -1.77%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D20' '1.004'
$ws.Range('E20').Value = '  +0.05%  '
Set-TextValue $ws 'D21' '4.425'
$ws.Range('E21').Value = '  -4.05%  '
$ws.Range('E22').Value = '  -3.15%  '
Set-TextValue $ws 'D23' '5.960'
$ws.Range('E23').Value = '  +0.13%  '
Set-TextValue $ws 'D24' '1.005'
$ws.Range('E24').Value = '  +0.08%  '
Set-TextValue $ws 'D25' '144.02'
$ws.Range('E25').Value = '  -0.13%  '
Set-TextValue $ws 'D26' '1.756'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('E28').Value = '  -3.27%  '
Set-TextValue $ws 'D29' '15.80'
$ws.Range('E29').Value = '  -1.71%  '
Set-TextValue $ws 'D30' '0.05258'
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('E31').Value = '  -2.55%  '
Set-TextValue $ws 'D32' '3.369'
$ws.Range('E32').Value = '  -2.72%  '
Set-TextValue $ws 'D33' '3.222'
$ws.Range('E33').Value = '  -5.22%  '
Set-TextValue $ws 'D34' '1.569'
$ws.Range('E34').Value = '  -4.46%  '
Set-TextValue $ws 'D35' '2.760'
$ws.Range('E35').Value = '  -3.59%  '
Set-TextValue $ws 'D36' '2.369'
$ws.Range('E36').Value = '  -1.91%  '
Set-TextValue $ws 'D37' '0.9277'
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('D38').Value = '1.164.85'
$ws.Range('E38').Value = '  +11.93%  '
Set-TextValue $ws 'D39' '0.5680'
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('E40').Value = '  -1.89%  '
Set-TextValue $ws 'D41' '0.8482'
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('E42').Value = '  +0.07%  '
Set-TextValue $ws 'D43' '5.689'
$ws.Range('E43').Value = '  -1.17%  '
Set-TextValue $ws 'D44' '100.52'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '1.797.16'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('E46').Value = '  -2.92%  '
Set-TextValue $ws 'D47' '0.4534'
$ws.Range('E47').Value = '  +0.45%  '
Set-TextValue $ws 'D48' '55.86'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('E49').Value = '  +0.26%  '
Set-TextValue $ws 'D50' '7.847'
$ws.Range('E50').Value = '  -2.42%  '
Set-TextValue $ws 'D51' '0.05053'
$ws.Range('E51').Value = '  -3.35%  '
